$d = $word.ActiveDocument

$replacements = @(
    @{old="39×53=2067"; new="14×68=952"},
    @{old="13×96=1248"; new="25×33=825"},
    @{old="34×95=3230"; new="28×94=2632"},
    @{old="16×58=928"; new="33×25=825"},
    @{old="64×55=3520"; new="40×27=1080"},
    @{old="23×95=2185"; new="89×68=6052"},
    @{old="48×71=3408"; new="90×27=2430"},
    @{old="23×22=506"; new="46×74=3404"},
    @{old="66×38=2508"; new="33×80=2640"},
    @{old="79×21=1659"; new="42×63=2646"},
    @{old="75×83=6225"; new="12×37=444"},
    @{old="12×33=396"; new="95×29=2755"},
    @{old="78×76=5928"; new="74×32=2368"},
    @{old="70×43=3010"; new="27×29=783"},
    @{old="72×63=4536"; new="27×34=918"},
    @{old="65×75=4875"; new="82×92=7544"},
    @{old="59×73=4307"; new="87×85=7395"},
    @{old="16×45=720"; new="34×99=3366"},
    @{old="30×65=1950"; new="74×89=6586"},
    @{old="92×55=5060"; new="72×84=6048"},
    @{old="63×58=3654"; new="81×59=4779"},
    @{old="84×82=6888"; new="99×53=5247"},
    @{old="41×64=2624"; new="56×57=3192"},
    @{old="12×53=636"; new="87×19=1653"},
    @{old="57×94=5358"; new="90×54=4860"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
